$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'257.09"
$ws.Range("E2").Value = "'-0.06%"
$ws.Range("E3").Value = "'-0.23%"
$ws.Range("D4").Value = "'4.684"
$ws.Range("E4").Value = "'-1.66%"
$ws.Range("E5").Value = "'-0.85%"
$ws.Range("E6").Value = "'-0.71%"
$ws.Range("D7").Value = "'0.8500"
$ws.Range("E7").Value = "'-2.57%"
$ws.Range("D8").Value = "'0.9280"
$ws.Range("E8").Value = "'-2.66%"
$ws.Range("D9").Value = "'0.1377"
$ws.Range("E9").Value = "'-2.04%"
$ws.Range("D10").Value = "'0.04271"
$ws.Range("E10").Value = "'9.45%"
$ws.Range("D11").Value = "'0.07031"
$ws.Range("E11").Value = "'-1.94%"
$ws.Range("D12").Value = "'0.03055"
$ws.Range("E12").Value = "'-4.44%"
$ws.Range("D13").Value = "'0.09106"
$ws.Range("E13").Value = "'-1.77%"
$ws.Range("D14").Value = "'0.001542"
$ws.Range("E14").Value = "'0.16%"
$ws.Range("D15").Value = "'0.0006068"
$ws.Range("E15").Value = "'0.13%"
$ws.Range("D16").Value = "'0.006013"
$ws.Range("E16").Value = "'-0.54%"
$ws.Range("E17").Value = "'-0.39%"
$ws.Range("D18").Value = "'3.177"
$ws.Range("E18").Value = "'-0.53%"
$ws.Range("D20").Value = "'0.3081"
$ws.Range("E20").Value = "'-1.71%"
$ws.Range("E21").Value = "'-0.93%"
$ws.Range("D22").Value = "'3.901"
$ws.Range("E22").Value = "'2.32%"
$ws.Range("D23").Value = "'0.04247"
$ws.Range("E23").Value = "'1.09%"
$ws.Range("D24").Value = "'0.001218"
$ws.Range("E24").Value = "'-0.49%"
$ws.Range("D25").Value = "'0.004298"
$ws.Range("E25").Value = "'-4.51%"
$ws.Range("D26").Value = "'0.0001200"
$ws.Range("E26").Value = "'-0.04%"
$ws.Range("E27").Value = "'2.00%"
$ws.Range("D40").Value = "'0.03794"
$ws.Range("E40").Value = "'-1.01%"
$ws.Range("D41").Value = "'0.006218"
$ws.Range("E41").Value = "'56.52%"
$ws.Range("E42").Value = "'-0.23%"
$ws.Range("E43").Value = "'-2.34%"
$ws.Range("E44").Value = "'32.68%"
$ws.Range("D45").Value = "'0.00005347"
$ws.Range("E45").Value = "'-2.78%"
$ws.Range("E46").Value = "'-0.03%"
$ws.Range("D47").Value = "'0.04314"
$ws.Range("E47").Value = "'-51.26%"
$ws.Range("E48").Value = "'10,465.23%"
$ws.Range("E49").Value = "'-0.03%"
$ws.Range("E50").Value = "'-0.03%"
